$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.909905365294291
$ws.Range("D2").Value = 10.01683104075481
$ws.Range("E2").Value = 9.229275407203918
$ws.Range("F2").Value = 73.44838837831928
$ws.Range("G2").Value = 3.826188401978788
$ws.Range("I2").Value = 20.31797844828558
$ws.Range("J2").Value = 8.569510664404456
$ws.Range("M2").Value = 41.77842852650505
$ws.Range("B3").Value = 7.598058114066005
$ws.Range("D3").Value = 9.667082402201361
$ws.Range("E3").Value = 8.807130319796354
$ws.Range("F3").Value = 72.93098895861142
$ws.Range("G3").Value = 3.83891309599214
$ws.Range("I3").Value = 20.28614089080616
$ws.Range("J3").Value = 8.590131279615292
$ws.Range("M3").Value = 40.74944359620653
$ws.Range("B4").Value = 7.399663309939654
$ws.Range("D4").Value = 9.452097644530205
$ws.Range("E4").Value = 8.536992819965842
$ws.Range("F4").Value = 72.65320759133024
$ws.Range("G4").Value = 3.847058570305241
$ws.Range("I4").Value = 20.27437417931135
$ws.Range("J4").Value = 8.60379077190299
$ws.Range("M4").Value = 40.11508236435017
$ws.Range("B5").Value = 7.317170824795721
$ws.Range("D5").Value = 9.364570211672708
$ws.Range("E5").Value = 8.424233855854981
$ws.Range("F5").Value = 72.54996009176428
$ws.Range("G5").Value = 3.85046245449664
$ws.Range("I5").Value = 20.27154272072922
$ws.Range("J5").Value = 8.609608663314674
$ws.Range("M5").Value = 39.85626642183419
$ws.Range("B6").Value = 7.303376668061976
$ws.Range("D6").Value = 9.350045485298383
$ws.Range("E6").Value = 8.40535115988191
$ws.Range("F6").Value = 72.53341405989356
$ws.Range("G6").Value = 3.851032799747114
$ws.Range("I6").Value = 20.27119131889746
$ws.Range("J6").Value = 8.610589926540712
$ws.Range("M6").Value = 39.81328164331502
$ws.Range("B7").Value = 7.398557313684842
$ws.Range("D7").Value = 9.450916700762845
$ws.Range("E7").Value = 8.535482838392456
$ws.Range("F7").Value = 72.65177499547242
$ws.Range("G7").Value = 3.847104132761294
$ws.Range("I7").Value = 20.2743280350839
$ws.Range("J7").Value = 8.603868214889976
$ws.Range("M7").Value = 40.11159267893073
$ws.Range("B8").Value = 7.803869831279236
$ws.Range("D8").Value = 9.896368437725469
$ws.Range("E8").Value = 9.086030983806735
$ws.Range("F8").Value = 73.26166320913671
$ws.Range("G8").Value = 3.830507446122104
$ws.Range("I8").Value = 20.3053903724683
$ws.Range("J8").Value = 8.576413741086549
$ws.Range("M8").Value = 41.42436272021644
$ws.Range("B9").Value = 8.540200595611285
$ws.Range("D9").Value = 10.7623446906824
$ws.Range("E9").Value = 10.07621868367634
$ws.Range("F9").Value = 74.77699230819042
$ws.Range("G9").Value = 3.800555662432168
$ws.Range("I9").Value = 20.42767179067203
$ws.Range("J9").Value = 8.530475247974499
$ws.Range("M9").Value = 43.96443646583738
$ws.Range("B10").Value = 9.041476246308873
$ws.Range("D10").Value = 11.38717643859981
$ws.Range("E10").Value = 10.74645156130895
$ws.Range("F10").Value = 76.08734563838119
$ws.Range("G10").Value = 3.780068970777322
$ws.Range("I10").Value = 20.55431101724097
$ws.Range("J10").Value = 8.501510533081998
$ws.Range("M10").Value = 45.79235496234704
$ws.Range("B11").Value = 9.260231164294817
$ws.Range("D11").Value = 11.66775337647631
$ws.Range("E11").Value = 11.03854174689617
$ws.Range("F11").Value = 76.72635863002752
$ws.Range("G11").Value = 3.771064703609857
$ws.Range("I11").Value = 20.61974798071489
$ws.Range("J11").Value = 8.489366815113392
$ws.Range("M11").Value = 46.61234145321343
$ws.Range("B12").Value = 9.341688986000195
$ws.Range("D12").Value = 11.77338996532677
$ws.Range("E12").Value = 11.14728592719599
$ws.Range("F12").Value = 76.97449082590209
$ws.Range("G12").Value = 3.767699163442143
$ws.Range("I12").Value = 20.64563648886609
$ws.Range("J12").Value = 8.484916286268209
$ws.Range("M12").Value = 46.92094572128372
$ws.Range("B13").Value = 9.324207612279857
$ws.Range("D13").Value = 11.7506676771809
$ws.Range("E13").Value = 11.12394912031973
$ws.Range("F13").Value = 76.92077799568241
$ws.Range("G13").Value = 3.768422046574107
$ws.Range("I13").Value = 20.64001189640749
$ws.Range("J13").Value = 8.48586821035995
$ws.Range("M13").Value = 46.85457076563903
$ws.Range("B14").Value = 9.266960640807666
$ws.Range("D14").Value = 11.67645686976184
$ws.Range("E14").Value = 11.04752566184182
$ws.Range("F14").Value = 76.7466494931453
$ws.Range("G14").Value = 3.770786939893621
$ws.Range("I14").Value = 20.62185571984087
$ws.Range("J14").Value = 8.488997703097247
$ws.Range("M14").Value = 46.63777006082643
$ws.Range("B15").Value = 9.231714275298684
$ws.Range("D15").Value = 11.63091860525985
$ws.Range("E15").Value = 11.00047076738604
$ws.Range("F15").Value = 76.64079104489568
$ws.Range("G15").Value = 3.772241224068836
$ws.Range("I15").Value = 20.61087843409198
$ws.Range("J15").Value = 8.490933873231848
$ws.Range("M15").Value = 46.5047181057106
$ws.Range("B16").Value = 9.026989325235894
$ws.Range("D16").Value = 11.3687588304662
$ws.Range("E16").Value = 10.72710303345354
$ws.Range("F16").Value = 76.04644697081449
$ws.Range("G16").Value = 3.78066366679663
$ws.Range("I16").Value = 20.55019068395504
$ws.Range("J16").Value = 8.50232489241383
$ws.Range("M16").Value = 45.73851301193049
$ws.Range("B17").Value = 8.898986017281576
$ws.Range("D17").Value = 11.20692994446076
$ws.Range("E17").Value = 10.55610314976116
$ws.Range("F17").Value = 75.69282046264291
$ws.Range("G17").Value = 3.78591049657307
$ws.Range("I17").Value = 20.51495426274336
$ws.Range("J17").Value = 8.509577045208484
$ws.Range("M17").Value = 45.26533165152985
$ws.Range("B18").Value = 8.824490866494969
$ws.Range("D18").Value = 11.11350988120479
$ws.Range("E18").Value = 10.45654444571183
$ws.Range("F18").Value = 75.49346657833713
$ws.Range("G18").Value = 3.788958076675661
$ws.Range("I18").Value = 20.49542512620544
$ws.Range("J18").Value = 8.513845499111182
$ws.Range("M18").Value = 44.99209711800892
$ws.Range("B19").Value = 8.799120001548689
$ws.Range("D19").Value = 11.08182378479496
$ws.Range("E19").Value = 10.42262959511793
$ws.Range("F19").Value = 75.42666351954628
$ws.Range("G19").Value = 3.789995077542555
$ws.Range("I19").Value = 20.48894012004116
$ws.Range("J19").Value = 8.51530743302377
$ws.Range("M19").Value = 44.89940818806087
$ws.Range("B20").Value = 8.912702690998749
$ws.Range("D20").Value = 11.22419278723376
$ws.Range("E20").Value = 10.57443118801242
$ws.Range("F20").Value = 75.73004630701323
$ws.Range("G20").Value = 3.785348892552209
$ws.Range("I20").Value = 20.51862898441274
$ws.Range("J20").Value = 8.508794983726027
$ws.Range("M20").Value = 45.31581554173415
$ws.Range("B21").Value = 9.28381324409038
$ws.Range("D21").Value = 11.69827160749406
$ws.Range("E21").Value = 11.07002382332346
$ws.Range("F21").Value = 76.79762853293678
$ws.Range("G21").Value = 3.770091123735945
$ws.Range("I21").Value = 20.6271586729643
$ws.Range("J21").Value = 8.488074481088537
$ws.Range("M21").Value = 46.70150328842399
$ws.Range("B22").Value = 9.518295850604426
$ws.Range("D22").Value = 12.0045121851769
$ws.Range("E22").Value = 11.38305124706941
$ws.Range("F22").Value = 77.53119373605011
$ws.Range("G22").Value = 3.760376241547795
$ws.Range("I22").Value = 20.70454414452278
$ws.Range("J22").Value = 8.475395087513844
$ws.Range("M22").Value = 47.59591757574357
$ws.Range("B23").Value = 9.393898841884893
$ws.Range("D23").Value = 11.84142009065855
$ws.Range("E23").Value = 11.2169831990876
$ws.Range("F23").Value = 77.13640740703737
$ws.Range("G23").Value = 3.765538137258025
$ws.Range("I23").Value = 20.66265732231737
$ws.Range("J23").Value = 8.482083528908618
$ws.Range("M23").Value = 47.11965347136166
$ws.Range("B24").Value = 8.906504196093291
$ws.Range("D24").Value = 11.21638944076128
$ws.Range("E24").Value = 10.56614896917921
$ws.Range("F24").Value = 75.71320419576453
$ws.Range("G24").Value = 3.785602696749953
$ws.Range("I24").Value = 20.51696537159708
$ws.Range("J24").Value = 8.509148245031534
$ws.Range("M24").Value = 45.29299548207846
$ws.Range("B25").Value = 8.34772148186965
$ws.Range("D25").Value = 10.52962995906025
$ws.Range("E25").Value = 9.818242802748349
$ws.Range("F25").Value = 74.33254656702114
$ws.Range("G25").Value = 3.808387033672609
$ws.Range("I25").Value = 20.38807804193976
$ws.Range("J25").Value = 8.542060180712912
$ws.Range("M25").Value = 43.28281804132899
